$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1,1).Value = "Índice"
$ws.Cells.Item(1,2).Value = "Distancia"
$ws.Cells.Item(1,3).Value = "max"
$ws.Cells.Item(1,4).Value = "min"
$ws.Cells.Item(1,5).Value = "Tempo"

# Data rows: Índice, Distancia, max, min, Tempo
$data = @(
    @(0, 2031, 2031, 2031, 0.01058477560679118),
    @(1, 2296, 2296, 2296, 0.01258331934611003),
    @(2, 1822, 1822, 1822, 0.01283709208170573),
    @(3, 2829, 2829, 2829, 0.01175345579783122),
    @(4, 2187, 2187, 2187, 0.01269187132517497),
    @(5, 3125, 3125, 3125, 0.01256766319274902),
    @(6, 2628, 2628, 2628, 0.01181619167327881),
    @(7, 2734, 2734, 2734, 0.01241788864135742),
    @(8, 2886, 2886, 2886, 0.01240902741750081),
    @(9, 2917, 2917, 2917, 0.01211783091227214)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
